# Commit: Fri, Jun 12, 2020  7:05:14 PM
#
# Changes applied:
#  1. Three tables (on slides 14, 15, 16) switch from the custom
#     "Table_0" table style to the built-in table style
#     {85E2D9D4-BA40-4531-8218-082E24BA75A2}.
#  2. The presentation's theme colour scheme (the "Integral" / Red Violet
#     theme used by the slide master) is swapped for the stock
#     "Office Theme" / Office colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newTableStyleId = "{85E2D9D4-BA40-4531-8218-082E24BA75A2}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the theme colour scheme ---------------------------------
# Office Theme / "Office" colour scheme values, in
# MsoThemeColorSchemeIndex order (Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink).
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [math]::Floor($hex / 65536) % 256
    $g = [math]::Floor($hex / 256) % 256
    $b = $hex % 256
    $bgr = ($b * 65536) + ($g * 256) + $r
    $themeColors.Item($i).RGB = $bgr
}
